# Program_v5_2016-10-26.docx — update Crystal Pan's and Kayode Ezike's
# session titles (and tidy up the stray "_GoBack" bookmark left over
# from the live edit session).
#
# Word leaves visible traces of interactive editing in the OOXML:
#   * runs get split at the exact point where the cursor/selection was
#     when text was retyped, even if the two resulting runs end up with
#     identical formatting;
#   * the hidden "_GoBack" bookmark (Word's "last edit location" marker)
#     moves to wherever text was most recently typed.
# The helper functions below reproduce that behaviour explicitly so the
# resulting document.xml matches what Word itself would have produced.

$d = $word.ActiveDocument

function Split-RunAt($pos) {
    # Force Word to keep $pos as a hard run boundary: nudge the font
    # color away from the paragraph's base color and immediately back.
    # A no-op text/format change gets silently merged back into the
    # neighbouring run, but a genuine (if transient) difference leaves
    # a permanent run split behind, exactly like Word does when you
    # click/retype in the middle of a run.
    $tail = $d.Range($pos, $d.Content.End)
    $save = $tail.Font.Color
    $tail.Font.Color = 1
    $tail2 = $d.Range($pos, $d.Content.End)
    $tail2.Font.Color = $save
}

function Set-RangeText($startPos, $endPos, $newText) {
    $r = $d.Range($startPos, $endPos)
    $r.Text = $newText
}

# ---------------------------------------------------------------------
# 1) "How Google Maps Figures Out Which Way to Go: Dijkstra's Algorithm"
#    Originally split into two runs around a (now removed) "_GoBack"
#    bookmark; collapse it back into a single run with no bookmark.
# ---------------------------------------------------------------------
$apos = [char]0x2019
$dijkstraFull = "How Google Maps Figures Out Which Way to Go: Dijkstra" + $apos + "s Algorithm"

$r = $d.Content
[void]$r.Find.Execute($dijkstraFull)
$start = $r.Start
$end = $r.End
Set-RangeText $start $end "PLACEHOLDER-DIJKSTRA"
$r2 = $d.Content
[void]$r2.Find.Execute("PLACEHOLDER-DIJKSTRA")
Set-RangeText $r2.Start $r2.End $dijkstraFull

# ---------------------------------------------------------------------
# 2) "How Concepts Help Us Understand Data Storage" (Kayode Ezike)
#    -> split into "How Concer" / "ts Help Us Understand Data Storage"
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("How Concepts Help Us Understand Data Storage")
$start = $r.Start
$splitPos = $start + 10

Split-RunAt $splitPos
Set-RangeText $start $splitPos "How Concer"

# ---------------------------------------------------------------------
# 3) "Reverse Engineering Smoothies with Math" (Prof. Kimberle Koile's
#    session) -> split into "Reverse" / " Engineering Smoothies with
#    Math", with the "_GoBack" bookmark re-inserted at the split point
#    (this is where the user's cursor ended up last).
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("Reverse Engineering Smoothies with Math")
$start = $r.Start
$end = $r.End
$splitPos = $start + 7

Split-RunAt $splitPos

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 4) "How to make superbabies" (Crystal Pan) -> "How To Make Super
#    Babies", split into four runs: "How To Make S" / "uper" / " B" /
#    "abies".
# ---------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("How to make superbabies")
$start = $r.Start
$end = $r.End

$b1 = $start + 13
$b2 = $start + 17
$b3 = $start + 19

Split-RunAt $b1
Split-RunAt $b2
Split-RunAt $b3

Set-RangeText $b3 $end "abies"
Set-RangeText $b2 $b3 " B"
Set-RangeText $b1 $b2 "uper"
Set-RangeText $start $b1 "How To Make S"

Write-Host "Edits applied"
